# [Outlook] (mapping) Include new snippets (#434)
# Append 19 new rows to the "Snippets" table describing additional
# Outlook Other-Item-APIs mapping entries (subject, internetMessageId,
# itemClass, itemType, start) across the Appointment/Message Read/Compose
# object models.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Columns: Class | Member Name | Member ID (methods only) | SnippetIdIntheYAMLFile | MethodNameInTheSnippet
$newSnippets = @(
    ,@("AppointmentRead", "subject", $null, "outlook-other-item-apis-get-subject-read", "get")
    ,@("MessageRead", "subject", $null, "outlook-other-item-apis-get-subject-read", "get")
    ,@("AppointmentCompose", "subject", $null, "outlook-other-item-apis-get-set-subject-compose", "get")
    ,@("MessageCompose", "subject", $null, "outlook-other-item-apis-get-set-subject-compose", "get")
    ,@("AppointmentCompose", "subject", $null, "outlook-other-item-apis-get-set-subject-compose", "set")
    ,@("MessageCompose", "subject", $null, "outlook-other-item-apis-get-set-subject-compose", "set")
    ,@("MessageRead", "internetMessageId", $null, "outlook-other-item-apis-get-internet-message-id-read", "get")
    ,@("AppointmentRead", "itemClass", $null, "outlook-other-item-apis-get-item-class-read", "get")
    ,@("MessageRead", "itemClass", $null, "outlook-other-item-apis-get-item-class-read", "get")
    ,@("AppointmentCompose", "itemType", $null, "outlook-other-item-apis-get-item-type", "get")
    ,@("AppointmentRead", "itemType", $null, "outlook-other-item-apis-get-item-type", "get")
    ,@("MessageCompose", "itemType", $null, "outlook-other-item-apis-get-item-type", "get")
    ,@("MessageRead", "itemType", $null, "outlook-other-item-apis-get-item-type", "get")
    ,@("AppointmentRead", "start", $null, "outlook-other-item-apis-get-start-read", "get")
    ,@("MessageRead", "start", $null, "outlook-other-item-apis-get-start-read", "get")
    ,@("AppointmentCompose", "start", $null, "outlook-other-item-apis-get-set-start-appointment-organizer", "get")
    ,@("Time", "getAsync", 2, "outlook-other-item-apis-get-set-start-appointment-organizer", "get")
    ,@("AppointmentCompose", "start", $null, "outlook-other-item-apis-get-set-start-appointment-organizer", "set")
    ,@("Time", "setAsync", 2, "outlook-other-item-apis-get-set-start-appointment-organizer", "set")
)

foreach ($entry in $newSnippets) {
    $newRow = $lo.ListRows.Add()
    $rng = $newRow.Range()
    $rng.Cells.Item(1, 1).Value() = $entry[0]
    $rng.Cells.Item(1, 2).Value() = $entry[1]
    if ($entry[2] -ne $null) {
        $rng.Cells.Item(1, 3).Value() = $entry[2]
    }
    $rng.Cells.Item(1, 4).Value() = $entry[3]
    $rng.Cells.Item(1, 5).Value() = $entry[4]
}

# Matches the authored file's resulting view state (frozen header row,
# top-left visible cell and selection reset to A2).
$ws.Range("A2").Select() | Out-Null